$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.497.76'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.578.36'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.000'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.42'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3692'
$ws.Range('E7').Value = '  +1.49%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '48.28'
$ws.Range('E8').Value = '  -3.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3340'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.153'
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07570'
$ws.Range('E11').Value = '  +2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.83'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.997'
$ws.Range('E14').Value = '  +1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.963'
$ws.Range('E15').Value = '  +1.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.581.26'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001127'
$ws.Range('E17').Value = '  +2.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '88.14'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06738'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.420'
$ws.Range('E20').Value = '  +1.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.67'
$ws.Range('E22').Value = '  +4.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.06'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '22.486.29'
$ws.Range('E24').Value = '  +0.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.399'
$ws.Range('E25').Value = '  +0.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.653'
$ws.Range('E26').Value = '  +4.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '152.12'
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.75'
$ws.Range('E28').Value = '  +0.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.993'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.64'
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.754.42'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.108'
$ws.Range('E32').Value = '  +4.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.180'
$ws.Range('E33').Value = '  +1.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.997'
$ws.Range('E34').Value = '  +1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.926'
$ws.Range('E35').Value = '  +4.01%  '
$ws.Range('E36').Value = '  +1.72%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02485'
$ws.Range('E37').Value = '  +4.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2264'
$ws.Range('E38').Value = '  +2.31%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06447'
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.411'
$ws.Range('E40').Value = '  +1.80%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.296'
$ws.Range('E41').Value = '  -1.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.54'
$ws.Range('E42').Value = '  +2.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.6328'
$ws.Range('E43').Value = '  +4.16%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.15'
$ws.Range('E44').Value = '  +2.79%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6163'
$ws.Range('E46').Value = '  +7.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.799'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.081'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '126.20'
$ws.Range('E49').Value = '  +1.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.220'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07250'
$ws.Range('E51').Value = '  +0.08%  '
